# Updates cryptocurrency price (column D) and 1h volume change (column E) values
# on the "cryptos" worksheet, matching the latest scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.994.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.790.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "358.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.228.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.817.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.945"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.927.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "274.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("E29").Value = "  +5.21%  "
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0466"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("E36").Value = "  +6.56%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.33%  "
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.076.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.939"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("E51").Value = "  +1.18%  "
